$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'88.213.53"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.56%  '

# Row 3
$ws.Range('D3').Value = "'3.266.97"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.18%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').Value = "'212.88"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.22%  '

# Row 6
$ws.Range('D6').Value = "'630.24"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.37%  '

# Row 7
$ws.Range('D7').Value = "'0.382"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +17.49%  '

# Row 8
$ws.Range('D8').Value = "'0.730"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +19.50%  '

# Row 9
$ws.Range('D9').Value = "'0.999"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '

# Row 10
$ws.Range('D10').Value = "'3.266.91"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.99%  '

# Row 11
$ws.Range('D11').Value = "'0.579"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.09%  '

# Row 12
$ws.Range('E12').Value = '  +11.92%  '

# Row 13
$ws.Range('D13').Value = "'0.0000264"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.03%  '

# Row 14
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D14').Value = "'5.53"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.27%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = "'34.47"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.75%  '

# Row 16
$ws.Range('D16').Value = "'3.865.62"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.14%  '

# Row 17
$ws.Range('D17').Value = "'87.794.69"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.17%  '

# Row 18
$ws.Range('D18').Value = "'3.294.46"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.63%  '

# Row 19
$ws.Range('D19').Value = "'3.16"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.11%  '

# Row 20
$ws.Range('E20').Value = '  -1.89%  '

# Row 21
$ws.Range('D21').Value = "'438.39"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.42%  '

# Row 22
$ws.Range('D22').Value = "'8.97"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.42%  '

# Row 23
$ws.Range('D23').Value = "'5.36"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.68%  '

# Row 24
$ws.Range('D24').Value = "'7.41"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.68%  '

# Row 25
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = "'5.28"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.78%  '

# Row 26
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').Value = "'12.34"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.43%  '

# Row 27
$ws.Range('D27').Value = "'3.451.56"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.60%  '

# Row 28
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = "'0.0000140"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.44%  '

# Row 29
$ws.Range('B29').Value = 'Litecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D29').Value = "'77.53"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.31%  '

# Row 30
$ws.Range('E30').Value = '  +0.10%  '

# Row 31
$ws.Range('D31').Value = "'0.179"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -16.70%  '

# Row 32
$ws.Range('E32').Value = '  -0.58%  '

# Row 33
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = "'8.86"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.35%  '

# Row 34
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = "'567.81"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.57%  '

# Row 35
$ws.Range('E35').Value = '  -8.86%  '

# Row 36
$ws.Range('D36').Value = "'7.19"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.53%  '

# Row 37
$ws.Range('E37').Value = '  -2.89%  '

# Row 38
$ws.Range('E38').Value = '  -6.53%  '

# Row 39
$ws.Range('D39').Value = "'22.96"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.83%  '

# Row 40
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').Value = "'21.82"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.96%  '

# Row 41
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.12%  '

# Row 42
$ws.Range('D42').Value = "'3.10"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.33%  '

# Row 43
$ws.Range('D43').Value = "'0.402"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.85%  '

# Row 44
$ws.Range('E44').Value = '  -0.31%  '

# Row 45
$ws.Range('E45').Value = '  +0.09%  '

# Row 46
$ws.Range('D46').Value = "'152.59"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.59%  '

# Row 47
$ws.Range('D47').Value = "'0.137"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +22.97%  '

# Row 48
$ws.Range('D48').Value = "'179.70"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.99%  '

# Row 49
$ws.Range('D49').Value = "'44.73"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.32%  '

# Row 50
$ws.Range('E50').Value = '  -1.33%  '

# Row 51
$ws.Range('D51').Value = "'4.25"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.03%  '
